# Update gh-pages output sheet data (南宁-漫展信息.xlsx) to the newly generated
# snapshot: bumps a handful of "want to go" counters, and adds the
# "南宁·小蜜蜂动漫嘉年华2.0" show ahead of the existing "南宁·AB动漫游戏嘉年华"
# / "横州·第二届海棠动漫游戏嘉年华" rows (which shift down by one).
#
# This touches two worksheets that both embed the same underlying exhibition
# rows: "展览" (Exhibition) and "全部类型" (All types, which interleaves the
# Exhibition/Performance/Local-life sheets sorted by date). The row offsets
# differ by one between the two sheets because "全部类型" carries an extra
# "演出" (Performance) row ahead of the insertion point.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibition) — insertion point is row 8
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

# Bump the "want to go" counters for the unaffected rows (2-7)
$ws.Range("F2").Value = 7120
$ws.Range("F3").Value = 58
$ws.Range("F4").Value = 203
$ws.Range("F5").Value = 125
$ws.Range("F6").Value = 1095
$ws.Range("F7").Value = 176

# Insert a new row at 8 — pushes the existing "南宁·AB动漫游戏嘉年华" (row 8)
# and "横州·第二届海棠动漫游戏嘉年华" (row 9) rows down to rows 9 and 10.
$ws.Rows.Item(8).Insert()

# The freshly inserted row loses the bordered/bold numbering style on column
# A; restore it by pasting the formatting (only) from the row above.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the new row 8 with the "南宁·小蜜蜂动漫嘉年华2.0" show. The "开始时间"
# column holds plain text like "2024-07-06", so force a text number format
# before assigning or Excel will silently reinterpret it as a date serial;
# clear the format again afterwards so the cell stays styled like its
# neighbours (no explicit number format applied).
$ws.Range("A8").Value = 7
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2024-07-06"
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = "南宁·小蜜蜂动漫嘉年华2.0"
$ws.Range("D8").Value = "亭洪路45号 百益上河城"
$ws.Range("E8").Value = "2024.07.06 10:00-07.06 17:00"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 50
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=84925"
$ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg"

# Row 9 now holds the shifted-down "南宁·AB动漫游戏嘉年华"; its "want to go"
# count grew from 31 to 48, and the running index in column A advances by
# one (it was 7, following the original row 8, and must become 8).
$ws.Range("F9").Value = 48
$ws.Range("A9").Value = 8

# Row 10 now holds the shifted-down "横州·第二届海棠动漫游戏嘉年华"; only its
# running index in column A advances (was 8, becomes 9).
$ws.Range("A10").Value = 9

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) — insertion point is row 9 (one row further
# down than "展览" because of the extra Performance row at row 8)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F2").Value = 7120
$ws.Range("F3").Value = 58
$ws.Range("F4").Value = 203
$ws.Range("F5").Value = 125
$ws.Range("F6").Value = 1095
$ws.Range("F7").Value = 176

$ws.Rows.Item(9).Insert()

$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A9").Value = 8
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2024-07-06"
$ws.Range("B9").ClearFormats()
$ws.Range("C9").Value = "南宁·小蜜蜂动漫嘉年华2.0"
$ws.Range("D9").Value = "亭洪路45号 百益上河城"
$ws.Range("E9").Value = "2024.07.06 10:00-07.06 17:00"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 50
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84925"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg"

$ws.Range("F10").Value = 48
$ws.Range("A10").Value = 9

$ws.Range("A11").Value = 10
